# Master_Template.xlsx - "Entity" sheet:
# Add a new "parent-catalog-item-id" line to the item-identification block,
# right before the "entity-information" section header (i.e. insert a new
# row 9, pushing the existing rows 9+ down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entity")

# xlShiftDown = -4121, xlPasteAll = -4104 (defaults used by Insert())
$xlShiftDown = -4121

# Insert a blank row above the current row 9 ("entity-information"),
# shifting rows 9-56 down to 10-57.
$ws.Rows.Item(9).Insert($xlShiftDown)

# The new row should look like the other "level 2" item-identification rows
# (e.g. row 8, "parent-catalog-item-type"), so copy that row's formatting
# down into the freshly inserted row 9.
$ws.Range("A8:D8").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new line item: level 2, tag "parent-catalog-item-id"; the
# value/instructions cells (C9/D9) stay blank, same as the template pattern.
$ws.Cells.Item(9, 1).Value = 2
$ws.Cells.Item(9, 2).Value = "parent-catalog-item-id"

# Leave the cursor on the cell that was just typed into.
$ws.Range("B9").Select()
